$wb = $excel.ActiveWorkbook
$betsWs = $wb.Worksheets.Item("bets")
$resWs  = $wb.Worksheets.Item("resumen")

# ---------------------------------------------------------------------------
# 1) "bets" sheet: append a new row 29 (bet #28) below the existing data.
# ---------------------------------------------------------------------------

$betsWs.Cells.Item(29, 1).Value = 28                    # A29 - bet id
$betsWs.Cells.Item(29, 2).Value = 45222                 # B29 - date
$betsWs.Cells.Item(29, 3).Value = 1                     # C29
$betsWs.Cells.Item(29, 4).Formula = "=F28"               # D29
$betsWs.Cells.Item(29, 5).Value = 230                   # E29
$betsWs.Cells.Item(29, 6).Formula = "=D29+E29"           # F29
$betsWs.Cells.Item(29, 7).Value = "ESPORTS"              # G29
$betsWs.Cells.Item(29, 8).Value = "WORLDS 2023"          # H29
$betsWs.Cells.Item(29, 9).Value = "DK"                   # I29 - new shared string
$betsWs.Cells.Item(29, 10).Value = "AMBOS TEAMS DRAKE"    # J29
$betsWs.Cells.Item(29, 11).Value = 1                     # K29
$betsWs.Cells.Item(29, 12).Value = 0                     # L29
$betsWs.Cells.Item(29, 13).Formula = "=ROUND((F29/`$D`$2-1)*100, 3)"  # M29

# Re-apply the date format to B29 and the numeric style to M29 (assigning a
# .Value/.Formula resets a cell's style, so formats have to be restored
# afterwards).
$betsWs.Range("B28").Copy($betsWs.Range("B29"))
$betsWs.Cells.Item(29, 2).Value = 45222
$betsWs.Range("M28").Copy($betsWs.Range("M29"))
$betsWs.Cells.Item(29, 13).Formula = "=ROUND((F29/`$D`$2-1)*100, 3)"

# ---------------------------------------------------------------------------
# 2) "resumen" sheet: E2 becomes a formula (+20000), and the array formulas
#    in B3:F3 need to be refreshed so they pick up the new last row of
#    bets!M (the freshly appended M29).
# ---------------------------------------------------------------------------

$resWs.Cells.Item(2, 5).Formula = "=12896.82+20000"

foreach ($col in 2..6) {
    $cell = $resWs.Cells.Item(3, $col)
    $f = $cell.Formula
    $cell.Formula = $f
}
